$wb = $excel.ActiveWorkbook
$wsFuel = $wb.Worksheets.Item("fuelprices")
$wsEmission = $wb.Worksheets.Item("emissionprices")

# --- Update existing fuel price values ---
$wsFuel.Range("D4").Value = 55
$wsFuel.Range("D14").Value = 75

# --- Append the new "Distributed Energy" / 2030 scenario rows to fuelprices ---
# (mirrors the existing "Distributed Energy" / 2040 block in rows 15-27)
$fuels  = @("Hardcoal","Nuclear","Gas","Heavyoil","Lightoil","Lignite","Biomass","Blackliquor","MSW","Oilshale","DRservice1","DRservice0","Hydrogencommod")
$prices = @(27, 2.5, 55, 65, 90, 9, 28, 1, 1, 10, 0, 0, 75)
# Columns whose cell keeps the "plain" (no explicit font color) look, by fuel index (0-based)
$plainCCol = @(9, 10, 11)   # Oilshale, DRservice1, DRservice0 -> C column default style
$plainDCol = @(12)          # Hydrogencommod -> D column default style

$row = 28
for ($i = 0; $i -lt $fuels.Length; $i++) {
    $cA = $wsFuel.Cells.Item($row, 1)
    $cB = $wsFuel.Cells.Item($row, 2)
    $cC = $wsFuel.Cells.Item($row, 3)
    $cD = $wsFuel.Cells.Item($row, 4)

    $cA.Value = "Distributed Energy"
    $cB.Value = 2030
    $cC.Value = $fuels[$i]
    $cD.Value = $prices[$i]

    $cA.Font.Color = 0
    $cB.Font.Color = 0
    if ($plainCCol -contains $i) {
        # leave C column at default font (no explicit color)
    } else {
        $cC.Font.Color = 0
    }
    if ($plainDCol -contains $i) {
        # leave D column at default font (no explicit color)
    } else {
        $cD.Font.Color = 0
    }

    $row++
}

# --- Append the new "Distributed Energy" / 2030 scenario row to emissionprices ---
# (mirrors the existing "Distributed Energy" / 2040 row, row 3)
$eA = $wsEmission.Cells.Item(4, 1)
$eB = $wsEmission.Cells.Item(4, 2)
$eC = $wsEmission.Cells.Item(4, 3)
$eD = $wsEmission.Cells.Item(4, 4)

$eA.Value = "Distributed Energy"
$eB.Value = 2030
$eC.Value = "CO2"
$eD.Value = 120

$eA.Font.Color = 0
$eB.Font.Color = 0
$eC.Font.Color = 0
# D4 keeps the default (no explicit color) font, matching D3's style

# --- Update sheet selections / active sheet ---
$wsFuel.Range("G11").Select()
$wsEmission.Select()
$wsEmission.Range("A13").Select()
